$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text so numeric-looking price strings (e.g. "0.999")
# are preserved as text, matching the source inlineStr cells, then strip
# the temporary number-format style so no stray style index is left behind.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range('D2').Value = '27.938.46'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.635.70'
$ws.Range('E3').Value = '  +0.07%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.15%  '
$ws.Range('D5').Value = '211.65'
$ws.Range('E5').Value = '  -0.20%  '
$ws.Range('E6').Value = '  -0.66%  '
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').Value = '1.867.15'
$ws.Range('E12').Value = '  +0.05%  '
$ws.Range('D13').Value = '1.644.81'
$ws.Range('E13').Value = '  +0.63%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.564'
$ws.Range('E14').Value = '  -0.42%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '4.01'
$ws.Range('E15').Value = '  -1.08%  '
$ws.Range('D16').Value = '65.29'
$ws.Range('E16').Value = '  -0.08%  '
$ws.Range('D17').Value = '27.950.47'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '229.67'
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('D19').Value = '7.82'
$ws.Range('E19').Value = '  +3.72%  '
$ws.Range('E20').Value = '  -0.42%  '
$ws.Range('E21').Value = '  -0.11%  '
$ws.Range('E22').Value = '  -0.27%  '
$ws.Range('D23').Value = '10.13'
$ws.Range('E23').Value = '  -2.36%  '
$ws.Range('E24').Value = '  +0.43%  '
$ws.Range('D25').Value = '156.17'
$ws.Range('E25').Value = '  +0.82%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').Value = '15.57'
$ws.Range('E28').Value = '  -0.57%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.31%  '
$ws.Range('D31').Value = '0.0481'
$ws.Range('E31').Value = '  -0.32%  '
$ws.Range('E32').Value = '  +0.69%  '
$ws.Range('E33').Value = '  +0.70%  '
$ws.Range('D34').Value = '1.401.91'
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('E35').Value = '  +2.86%  '
$ws.Range('E36').Value = '  +0.50%  '
$ws.Range('E37').Value = '  -0.85%  '
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('E39').Value = '  -0.48%  '
$ws.Range('D40').Value = '0.853'
$ws.Range('E40').Value = '  -2.30%  '
$ws.Range('E41').Value = '  -0.09%  '
$ws.Range('D42').Value = '1.01'
$ws.Range('E42').Value = '  -0.95%  '
$ws.Range('E43').Value = '  +1.46%  '
$ws.Range('D44').Value = '66.08'
$ws.Range('E44').Value = '  -1.40%  '
$ws.Range('E45').Value = '  -1.74%  '
$ws.Range('D46').Value = '1.775.34'
$ws.Range('E46').Value = '  -0.08%  '
$ws.Range('E47').Value = '  -2.63%  '
$ws.Range('D48').Value = '88.58'
$ws.Range('E48').Value = '  +0.79%  '
$ws.Range('E49').Value = '  +2.08%  '
$ws.Range('D50').Value = '0.0506'
$ws.Range('E50').Value = '  -0.03%  '
$ws.Range('E51').Value = '  +1.21%  '

$colD.ClearFormats()
